$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoWebShop")

# New "Order number" column (F) with header + value for the new test case row
$ws.Range("F1").Value = "Order number"
$ws.Range("F3").Value = "'997040"

# Match the bestFit-style width Excel computed for the new column
$ws.Columns.Item(6).ColumnWidth = 11.5

# Move the active selection as recorded in the saved workbook
$ws.Range("C12").Select()
